# Add new "html" file-format master-data rows (eng/ara/fra) to the
# file_format master template sheet, mirroring the existing txt/xml/json
# rows (code, descr, lang_code, is_active, cr_by, cr_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: 11 = eng, 12 = ara, 13 = fra.
# Column-by-column entry order (code for all 3 rows, then descr for ara
# before eng, etc.) reproduces the exact shared-string intern order the
# author's edit produced.
$ws.Cells.Item(11, 1).Value = "html"
$ws.Cells.Item(12, 1).Value = "html"
$ws.Cells.Item(13, 1).Value = "html"

$ws.Cells.Item(12, 2).Value = "ملف html"
$ws.Cells.Item(11, 2).Value = "html file"
$ws.Cells.Item(13, 2).Value = "Fichier html"

$ws.Cells.Item(11, 3).Value = "eng"
$ws.Cells.Item(12, 3).Value = "ara"
$ws.Cells.Item(13, 3).Value = "fra"

$ws.Cells.Item(11, 4).Value = $true
$ws.Cells.Item(12, 4).Value = $true
$ws.Cells.Item(13, 4).Value = $true

$ws.Cells.Item(11, 5).Value = "superadmin"
$ws.Cells.Item(12, 5).Value = "superadmin"
$ws.Cells.Item(13, 5).Value = "superadmin"

$ws.Cells.Item(11, 6).Value = "now()"
$ws.Cells.Item(12, 6).Value = "now()"
$ws.Cells.Item(13, 6).Value = "now()"

# Match the existing left-aligned style used on the is_active column.
$ws.Range("D11:D13").HorizontalAlignment = -4131

# Mirror the author's final selection (columns G through the last column,
# XFD) as seen in the diff.
$ws.Range("G1:XFD1048576").Select()
